$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, pushing the existing row 86 (and everything
# below it, through the former row 154) down by one -- the former row 154
# ends up landing in row 155, growing the used range to A1:R155.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new weekly record.
$ws.Cells.Item(86, 1).Value = 11
$ws.Cells.Item(86, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(86, 3).Value = "Bíobío"
$ws.Cells.Item(86, 4).Value = 44484
$ws.Cells.Item(86, 5).Value = 8
$ws.Cells.Item(86, 6).Value = 100112008
$ws.Cells.Item(86, 7).Value = "Coliflor"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 1500
$ws.Cells.Item(86, 11).Value = 600
$ws.Cells.Item(86, 12).Value = 700
$ws.Cells.Item(86, 13).Value = 653
$ws.Cells.Item(86, 14).Value = "$/unidad"
$ws.Cells.Item(86, 15).Value = "Región Metropolitana"
$ws.Cells.Item(86, 16).Value = 653
$ws.Cells.Item(86, 17).Value = 1
$ws.Cells.Item(86, 18).Value = "Hortaliza"
